$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("E1").Value = "Generate Remittances"
$ws.Range("F1").Value = "Remittances Verified"

# New data cells
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "Yes"

$ws.Range("E3").Value = "Yes"
$ws.Range("F3").Value = "No"

$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "No"

# Column widths to match bestFit sizing (closest achievable given width quantization)
$ws.Columns.Item(5).ColumnWidth = 18.285714285714285
$ws.Columns.Item(6).ColumnWidth = 17

# Update selection to match target state
$ws.Range("E5").Select()
